# Apply weekly update to the "Apio" (Hortaliza) sheet.
# New observations for week of 2021-10-xx get inserted right after the
# existing "Provincia del Elquí" (44246) pair, pushing the rest of the
# historical data down by 4 rows (2 new date-pairs inserted).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows starting at row 365; this shifts the existing
# rows 365:425 down to 369:429 (and carries the D-column date style
# down with them automatically).
$ws.Rows("365:368").Insert()

# Fill in the 4 newly inserted rows with the new data.
$newRows = @(
    @{ Row=365; A=6; B="Mercado Mayorista Lo Valledor de Santiago"; C="Metropolitana"; D=44505; E=13; F=100112017; G="Apio"; H="Americana (o)"; I="Primera"; J=690;  K=5000; L=6000; M=5536; N="`$/docena de matas"; O="Región Metropolitana"; P=923; Q=6; R="Hortaliza" },
    @{ Row=366; A=6; B="Mercado Mayorista Lo Valledor de Santiago"; C="Metropolitana"; D=44505; E=13; F=100112017; G="Apio"; H="Americana (o)"; I="Primera"; J=2800; K=5000; L=6000; M=5571; N="`$/docena de matas"; O="Región de Coquimbo"; P=928; Q=6; R="Hortaliza" },
    @{ Row=367; A=6; B="Mercado Mayorista Lo Valledor de Santiago"; C="Metropolitana"; D=44505; E=13; F=100112017; G="Apio"; H="Americana (o)"; I="Segunda"; J=260;  K=4000; L=4000; M=4000; N="`$/docena de matas"; O="Región Metropolitana"; P=667; Q=6; R="Hortaliza" },
    @{ Row=368; A=6; B="Mercado Mayorista Lo Valledor de Santiago"; C="Metropolitana"; D=44505; E=13; F=100112017; G="Apio"; H="Americana (o)"; I="Segunda"; J=800;  K=4000; L=4000; M=4000; N="`$/docena de matas"; O="Región de Coquimbo"; P=667; Q=6; R="Hortaliza" }
)

foreach ($rowData in $newRows) {
    $r = $rowData.Row
    $arr = New-Object 'object[,]' 1,18
    $arr[0,0]  = $rowData.A
    $arr[0,1]  = $rowData.B
    $arr[0,2]  = $rowData.C
    $arr[0,3]  = $rowData.D
    $arr[0,4]  = $rowData.E
    $arr[0,5]  = $rowData.F
    $arr[0,6]  = $rowData.G
    $arr[0,7]  = $rowData.H
    $arr[0,8]  = $rowData.I
    $arr[0,9]  = $rowData.J
    $arr[0,10] = $rowData.K
    $arr[0,11] = $rowData.L
    $arr[0,12] = $rowData.M
    $arr[0,13] = $rowData.N
    $arr[0,14] = $rowData.O
    $arr[0,15] = $rowData.P
    $arr[0,16] = $rowData.Q
    $arr[0,17] = $rowData.R
    $ws.Range("A$r`:R$r").Value2 = $arr
}

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()
